# Nexial unitTest_web.xlsx update:
#  - add new JSON command `storeKeys(json,jsonpath,var)` to the `#system`
#    sheet's `json` list (alphabetically before `storeValue`)
#  - retire the single-member `text` list (it only ever held
#    `spellCheck(var,profile,text)`), which shifts every column from the
#    old `web` list onward one column to the left
#  - keep the `target` meta-list (the list of all list-names) in sync by
#    dropping its `text` entry too
#  - update the `definedNames` ranges that describe all of this

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Drop the old column Y ("text"). Excel's column-delete semantics
#    shift Z->Y, AA->Z, AB->AA, AC->AB, AD->AC, AE->AD and clear the
#    vacated column AE -- exactly the column-family shift that the
#    `web`/`webalert`/`webcookie`/`ws`/`ws.async`/`xml` lists need.
# ---------------------------------------------------------------------
$ws.Columns("Y:Y").Delete()

# ---------------------------------------------------------------------
# 2) Insert the new `storeKeys(json,jsonpath,var)` entry into the `json`
#    list (column M), directly above `storeValue`, pushing the rest of
#    the (alphabetically sorted) list down by one row.
# ---------------------------------------------------------------------
$ws.Range("M18").Value = "storeValues(json,jsonpath,var)"
$ws.Range("M17").Value = "storeValue(json,jsonpath,var)"
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------
# 3) The `target` column (A) mirrors every list's own name; remove its
#    "text" entry (row 25) and shift the remaining names up by one row.
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "web"
$ws.Range("A26").Value = "webalert"
$ws.Range("A27").Value = "webcookie"
$ws.Range("A28").Value = "ws"
$ws.Range("A29").Value = "ws.async"
$ws.Range("A30").Value = "xml"
$ws.Range("A31").ClearContents()

# ---------------------------------------------------------------------
# 4) Sync up the workbook-level named ranges.
# ---------------------------------------------------------------------
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
